# Update the EC (Estado de Cuenta) database values for NIT-8909401641.
# The worksheet lists overdue-period rows (column E = period code, column F = "Valor Mora").
# Row 16 corresponds to period 1807 and row 30 corresponds to period 1705; their
# "Valor Mora" amounts are swapped as part of the database refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Period 1807 (row 16): Valor Mora 15700 -> 31400
$ws.Range("F16").Value = 31400

# Period 1705 (row 30): Valor Mora 31400 -> 15700
$ws.Range("F30").Value = 15700
